$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A: Mat (student id numbers)
$ws.Range("A2").Value = 19330051920085
$ws.Range("A3").Value = 19330051920093
$ws.Range("A4").Value = 19330051920097
$ws.Range("A5").Value = 19330051920103
$ws.Range("A6").Value = 19330051920110
$ws.Range("A7").Value = 19330051920116
$ws.Range("A8").Value = 19330051920121
$ws.Range("A9").Value = 19330051920122
$ws.Range("A10").Value = 19330051920084
$ws.Range("A11").Value = 19330051920088
$ws.Range("A12").Value = 19330051920090
$ws.Range("A13").Value = 19330051920089
$ws.Range("A14").Value = 19330051920100
$ws.Range("A15").Value = 19330051920106
$ws.Range("A16").Value = 19330051920107
$ws.Range("A17").Value = 19330051920115

# Column B: Paterno
$ws.Range("B2").Value = "ANTONIO"
$ws.Range("B3").Value = "CONDADO"
$ws.Range("B4").Value = "FLORENCIO"
$ws.Range("B5").Value = "HERNANDEZ"
$ws.Range("B6").Value = "MARTINEZ"
$ws.Range("B7").Value = "ROJAS"
$ws.Range("B8").Value = "TEXCAHUA"
$ws.Range("B9").Value = "TOCOHUA"
$ws.Range("B10").Value = "AGUILAR"
$ws.Range("B11").Value = "CADEZA"
$ws.Range("B12").Value = "CASTILLO"
$ws.Range("B13").Value = "CASTILLO"
$ws.Range("B14").Value = "FLORES"
$ws.Range("B15").Value = "JIMENEZ"
$ws.Range("B16").Value = "JUAN"
$ws.Range("B17").Value = "OCTAVIANO"

# Column C: Materno
$ws.Range("C2").Value = "TEXOCO"
$ws.Range("C3").Value = "MORALES"
$ws.Range("C4").Value = "BERNABE"
$ws.Range("C5").Value = "DEL ANGEL"
$ws.Range("C6").Value = "RODRIGUEZ"
$ws.Range("C7").Value = "ANGUIANO"
$ws.Range("C8").Value = "TETLA"
$ws.Range("C9").Value = "BERISTAIN"
$ws.Range("C10").Value = "TEZOCO"
$ws.Range("C11").Value = "GALLARDO"
$ws.Range("C12").Value = "ROJAS"
$ws.Range("C13").Value = "ROJAS"
$ws.Range("C14").Value = "SANCHEZ"
$ws.Range("C15").Value = "MIRON"
$ws.Range("C16").Value = "DE LA CRUZ"
$ws.Range("C17").Value = "HERNANDEZ"

# Column D: Nombres
$ws.Range("D2").Value = "JOSE JAZAEL"
$ws.Range("D3").Value = "JOSUE ALEXIS"
$ws.Range("D4").Value = "JESUS ENRIQUE"
$ws.Range("D5").Value = "URIEL"
$ws.Range("D6").Value = "YAEL"
$ws.Range("D7").Value = "JOSUE"
$ws.Range("D8").Value = "JOSUE GILBERTO"
$ws.Range("D9").Value = "ADOLFO"
$ws.Range("D10").Value = "JOSE ANTONIO"
$ws.Range("D11").Value = "IRVING OTTONIEL"
$ws.Range("D12").Value = "BRAYAN"
$ws.Range("D13").Value = "BRANDON"
$ws.Range("D14").Value = "JUAN LUIS"
$ws.Range("D15").Value = "EMMANUEL"
$ws.Range("D16").Value = "OSWALDO ENRIQUE"
$ws.Range("D17").Value = "JOSE ANTONIO"

# Column E: Nombre_Largo (materia)
$ws.Range("E2").Value = "INGLÉS IV"
$ws.Range("E3").Value = "INGLÉS IV"
$ws.Range("E4").Value = "INGLÉS IV"
$ws.Range("E5").Value = "INGLÉS IV"
$ws.Range("E6").Value = "INGLÉS IV"
$ws.Range("E7").Value = "INGLÉS IV"
$ws.Range("E8").Value = "INGLÉS IV"
$ws.Range("E9").Value = "INGLÉS IV"
$ws.Range("E10").Value = "INGLÉS IV"
$ws.Range("E11").Value = "INGLÉS IV"
$ws.Range("E12").Value = "INGLÉS IV"
$ws.Range("E13").Value = "INGLÉS IV"
$ws.Range("E14").Value = "INGLÉS IV"
$ws.Range("E15").Value = "INGLÉS IV"
$ws.Range("E16").Value = "INGLÉS IV"
$ws.Range("E17").Value = "INGLÉS IV"

# Column F: Grupo
$ws.Range("F2").Value = "4BEM"
$ws.Range("F3").Value = "4BEM"
$ws.Range("F4").Value = "4BEM"
$ws.Range("F5").Value = "4BEM"
$ws.Range("F6").Value = "4BEM"
$ws.Range("F7").Value = "4BEM"
$ws.Range("F8").Value = "4BEM"
$ws.Range("F9").Value = "4BEM"
$ws.Range("F10").Value = "4BEM"
$ws.Range("F11").Value = "4BEM"
$ws.Range("F12").Value = "4BEM"
$ws.Range("F13").Value = "4BEM"
$ws.Range("F14").Value = "4BEM"
$ws.Range("F15").Value = "4BEM"
$ws.Range("F16").Value = "4BEM"
$ws.Range("F17").Value = "4BEM"

# Column G: Reprobadas (count)
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("G17").Value = 1

